$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "60.356.00"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "2.602.85"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.607.21"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.156"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.372"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "3.055.78"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("D16").Value = "60.333.05"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "2.603.26"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.97%  "
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +4.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.520"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.98%  "
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("E30").Value = "  +9.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.40%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "164.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "310.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.840"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.606"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0549"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0242"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.67%  "
